# Add a "year" column (column C) to each of the per-year worksheets
# (sheets named "2009".."2024"), filling it with that sheet's year for
# every data row. The Summary sheet already has a "year" column and is
# left untouched.

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Only touch the per-year data sheets (named e.g. "2009"); skip the
    # "Summary" sheet (and anything else whose name isn't a bare year).
    if ($ws.Name -notmatch '^[0-9]+$') {
        continue
    }
    $year = [int]$ws.Name

    $lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1

    # Header cell: label + same formatting as the existing header cells.
    $ws.Cells.Item(1, 1).Copy() | Out-Null
    $ws.Cells.Item(1, 3).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item(1, 3).Value = "year"

    # Data cells: numeric year value, one per existing data row.
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 3).Value = $year
    }
}

$excel.CutCopyMode = 0
